$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

# --- zh-cn sheet ---
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("B3").Value = "Ready for handoff"
$zh.Range("C3").Value = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.zh-cn.xlf"
$zh.Range("D3").Value = "2016-01-25 06:47:51"

# update hyperlink display text for C3 on zh-cn (keep same target)
$zhLink = $zh.Hyperlinks.Item(1)
foreach ($hl in $zh.Hyperlinks) {
    if ($hl.Range.Address -eq "$C3") {
        $hl.TextToDisplay = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.zh-cn.xlf"
    }
}

# --- de-de sheet ---
$de = $wb.Worksheets.Item("de-de")
$de.Range("B3").Value = "Ready for handoff"
$de.Range("C3").Value = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.de-de.xlf"
$de.Range("D3").Value = "2016-01-25 06:48:02"

foreach ($hl in $de.Hyperlinks) {
    if ($hl.Range.Address -eq "$C3") {
        $hl.TextToDisplay = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.de-de.xlf"
    }
}
